# Update "想去人数" (number of people interested) counts that changed
# between the previous generation and this one (456a3b4).
#
# Sheet "展览" (Exhibitions) and sheet "全部类型" (All types) both list the
# same two exhibitions; their F2/F3 cells need to be bumped:
#   F2: 149 -> 152
#   F3: 41  -> 43

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 152
    $ws.Range("F3").Value = 43
}
